# Add "Tasks Page" section (with its two sub-headings / bullet list) right
# after the "Notification Settings Page" section's last bullet, and link
# the Heading1/Heading2 paragraph styles to new "*Char" character styles.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the new "Tasks Page" block of paragraphs.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Add notification settings as panels are updated", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorPara = $rng.Paragraphs(1)
$pageBreakPara = $anchorPara.Next()
$insertionPoint = $d.Range($pageBreakPara.Range.Start, $pageBreakPara.Range.Start)

# Raw OOXML fragment for the new paragraphs: two blank paragraphs, the
# "Tasks Page" / "Adding new task" headings, three bulleted list items and
# a trailing blank paragraph. Using InsertXML (rather than typing "`r")
# keeps the untouched page-break paragraph that follows completely intact
# and avoids leaving stray empty runs behind in the blank paragraphs.
$xmlFrag = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:t>Tasks Page</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t>Adding new task</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr><w:r><w:t>And have the contect person that’s linked to that item for the task auto selected based on the linked project and have a add icon to add more people for the task</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr><w:r><w:t>Make all dropdown fields even</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr><w:r><w:t>Need to make a lot more changes check pipdrive for items that’s needed</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
'@

$insertionPoint.InsertXML($xmlFrag) | Out-Null

# ---------------------------------------------------------------------
# 2) Link Heading1 / Heading2 to new linked character styles and define
#    those character styles (mirrors the paragraph styles' run formatting).
# ---------------------------------------------------------------------
$heading1 = $d.Styles("Heading1")
$heading2 = $d.Styles("Heading2")

$heading1Char = $d.Styles.Add("Heading1Char", 2)
$heading1Char.NameLocal = "Heading 1 Char"
$heading1Char.BaseStyle = "DefaultParagraphFont"
$heading1Char.LinkStyle = "Heading1"
$heading1Char.Priority = 9
$heading1Char.Font.Bold = $true
$heading1Char.Font.BoldBi = $true
$heading1Char.Font.Color = 11485214
$heading1Char.Font.Size = 18
$heading1Char.Font.SizeBi = 18

$heading2Char = $d.Styles.Add("Heading2Char", 2)
$heading2Char.NameLocal = "Heading 2 Char"
$heading2Char.BaseStyle = "DefaultParagraphFont"
$heading2Char.LinkStyle = "Heading2"
$heading2Char.Priority = 9
$heading2Char.Font.Bold = $true
$heading2Char.Font.BoldBi = $true
$heading2Char.Font.Color = 5587251
$heading2Char.Font.Size = 14
$heading2Char.Font.SizeBi = 14

$heading1.LinkStyle = "Heading1Char"
$heading2.LinkStyle = "Heading2Char"

Write-Output "Tasks Page section + linked heading character styles added."
